$wb = $excel.ActiveWorkbook

# --- Sheet1 updates ---
$ws1 = $wb.Worksheets.Item("Sheet1")

$sheet1Data = @(
    @(1, 6),
    @(3, 1),
    @(6, 1),
    @(55, 1),
    @(2145, 1),
    @(4186, 2),
    @(4560, 1),
    @(4950, 1),
    @(5151, 1),
    @(134421, 1)
)

$r = 2
foreach ($row in $sheet1Data) {
    $ws1.Cells.Item($r, 1).Value = $row[0]
    $ws1.Cells.Item($r, 2).Value = $row[1]
    $r = $r + 1
}

# --- Sheet2 updates ---
$ws2 = $wb.Worksheets.Item("Sheet2")

$sheet2Data = @(
    @(1, 5),
    @(3, 1),
    @(15, 2),
    @(55, 1),
    @(91, 1),
    @(2701, 1),
    @(5565, 1),
    @(361675, 1)
)

$r = 2
foreach ($row in $sheet2Data) {
    $ws2.Cells.Item($r, 1).Value = $row[0]
    $ws2.Cells.Item($r, 2).Value = $row[1]
    $r = $r + 1
}

# Remove now-unused rows 10-15 (previous data went down to row 15)
$ws2.Range("A10:B15").ClearContents() | Out-Null
